$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide content: "JUnit, Qunit, Selenium" -> "TestNG, QUnit, Selenium"
#    (Issue 769: Migrate from JUnit to TestNG; also fix "Qunit" -> "QUnit")
#    Edit the individual runs in place (via Characters) so that only the
#    targeted run text changes and all run/paragraph formatting is preserved.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(17)
$tr = $sh.TextFrame.TextRange

$tr.Characters(1, 5).Text = "TestNG"
$tr.Characters(9, 5).Text = "QUnit"

# ---------------------------------------------------------------------------
# 2) Footer "date" placeholder text on the slide master and every slide
#    layout: "3/24/2013" -> "4/28/2013"
#    (Note: going through the slide master reached via Presentation.Designs
#    is required here so each CustomLayouts.Item(i) resolves to its own,
#    distinct layout.)
# ---------------------------------------------------------------------------
$sm = $p.Designs.Item(1).SlideMaster

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shape = $shapes.Item($k)
        if ($shape.Name -like "Date Placeholder*" -and $shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "3/24/2013") {
                $shape.TextFrame.TextRange.Text = "4/28/2013"
            }
        }
    }
}

Update-DatePlaceholder $sm.Shapes

$layouts = $sm.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}
